$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Recomputed line loading percentages for the 380 kV case (Case_2_253).
# Data rows correspond to A2:A25 (time steps 0..23); only columns B, C, E, F, G, I, L
# change between the old and new run. Columns D, H, J, K, M, N, O remain 0.

$newValues = @{
    "B" = @(21.37456280765034, 21.02796891469282, 20.82044561178383, 20.7373306767611, 20.72362077553542, 20.81931865290227, 21.25403456319298, 22.14295455326275, 22.81085997330377, 23.11638019881017, 23.23219052715968, 23.2072454834453, 23.12590635739821, 23.07609528963361, 22.79091657497641, 22.61630967045607, 22.51604646802305, 22.48213150249317, 22.63488061255941, 23.14979545361012, 23.48694552773062, 23.3069856540483, 22.6264843042762, 21.89934381485559)
    "C" = @(11.36548476683183, 10.8147708032762, 10.46587687577909, 10.32121285054541, 10.29704826769228, 10.46393564529525, 11.17795106704996, 12.48479792636159, 13.37919671830118, 13.77040480060559, 13.9162021014458, 13.88490763566275, 13.78244709906592, 13.71937922837659, 13.35330692980091, 13.12465138486542, 12.99166533963889, 12.94638928395323, 13.14914496402416, 13.81260653823462, 14.23251809874811, 14.00968349634848, 13.13807617192364, 12.14226127621667)
    "E" = @(8.625279463289823, 8.63470932158665, 8.640936578185881, 8.643584358786596, 8.644030676824569, 8.640971840972119, 8.628440217106363, 8.607328340265493, 8.593919066104352, 8.588273350148519, 8.586200650951332, 8.586644145171366, 8.588101521599805, 8.589002696080792, 8.594297157903178, 8.597661399790017, 8.599639182995565, 8.600316173894914, 8.597298845770322, 8.587671685828917, 8.581759834776193, 8.584880359420419, 8.597462620518662, 8.612669949220047)
    "F" = @(16.86991607391233, 15.89584955866808, 15.26997757108491, 15.00819731993403, 14.96433081551593, 15.26647399323728, 16.53996406344769, 19.0027458068253, 20.67494806633232, 21.3917225636224, 21.65686569030329, 21.60004134736742, 21.4136618050453, 21.29868154950795, 20.62722412089977, 20.20408069617459, 19.95656407809808, 19.87204792380562, 20.2495528364879, 21.46857628470567, 22.22866616901555, 21.82633154475864, 20.22900810905294, 18.34778573295697)
    "G" = @(3.716552107095005, 3.721386874987418, 3.724502955418774, 3.725810044829408, 3.726029341439807, 3.724520432190013, 3.718188620060878, 3.706934731895992, 3.69936458713476, 3.696069998064538, 3.694843684039846, 3.69510684912612, 3.695968683034555, 3.696499346931402, 3.699582882508327, 3.701512603610348, 3.702636572825582, 3.703019546732666, 3.701305729137731, 3.695714965442256, 3.692185009945063, 3.694057728816777, 3.701399211764071, 3.70985583693757)
    "I" = @(35.83496705610499, 35.48134649904628, 35.26859887381717, 35.18305344719466, 35.16891962016688, 35.26744045381419, 35.71216138659544, 36.61664623563702, 37.29763176794916, 37.61026718639508, 37.72899999233474, 37.70341448346583, 37.62002887783107, 37.56899570291958, 37.27725210746235, 37.09896087661425, 36.99668487073117, 36.96210474041601, 37.11791256306913, 37.64451242550146, 37.9906534819969, 37.80575232027429, 37.10934379673004, 36.3688225646951)
    "L" = @(10.77073310614577, 10.75407328103996, 10.74624333475613, 10.74365683372795, 10.7432638621123, 10.74620600471694, 10.76449153952182, 10.81932598173085, 10.87107372354665, 10.89707120294673, 10.90726555438888, 10.90505452990502, 10.89790291306917, 10.89356776210511, 10.8694238684667, 10.8552392465673, 10.8473121566668, 10.8446680499673, 10.85672528650825, 10.89999405780082, 10.93030883188128, 10.91394428730065, 10.85605273791296, 10.80246948647728)
}

$startRow = 2
foreach ($col in $newValues.Keys) {
    $colIndex = $ws.Range($col + "1").Column
    $values = $newValues[$col]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($startRow + $i, $colIndex).Value = $values[$i]
    }
}
